$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 400
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 400
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 400
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -740

$ws.Range("H17").Value = 985.08826
$ws.Range("I17").Value = 950.2353000000001
$ws.Range("J17").Value = 1019.94116
$ws.Range("K17").Value = 2850.7059
$ws.Range("L17").Value = 3059.82348
$ws.Range("M17").Value = -2682.7059
$ws.Range("N17").Value = -3395.82348

$ws.Range("H39").Value = 23.09091
$ws.Range("I39").Value = 23.09091
$ws.Range("J39").Value = 0
$ws.Range("K39").Value = 69.27273
$ws.Range("L39").Value = 0
$ws.Range("M39").Value = 226.72727
$ws.Range("N39").ClearContents()

$ws.Range("H43").Value = 5899.9414
$ws.Range("I43").Value = 4787.375
$ws.Range("J43").Value = 6888.8887
$ws.Range("K43").Value = 4787.375
$ws.Range("L43").Value = 6888.8887
$ws.Range("M43").Value = -4718.375
$ws.Range("N43").Value = -7026.8887

$ws.Range("H70").Value = 176354.5
$ws.Range("J70").Value = 234139.33
$ws.Range("L70").Value = 702417.99
$ws.Range("N70").Value = -702957.99

$ws.Range("H73").Value = 176354.5
$ws.Range("J73").Value = 234139.33
$ws.Range("L73").Value = 702417.99
$ws.Range("N73").Value = -704289.99

$ws.Range("H98").Value = 1290.5927
$ws.Range("I98").Value = 1313.88
$ws.Range("K98").Value = 1313.88
$ws.Range("M98").Value = 184.1199999999999

$ws.Range("H99").Value = 381.66666
$ws.Range("I99").Value = 381.66666
$ws.Range("K99").Value = 1144.99998
$ws.Range("M99").Value = 353.0000199999999

$ws.Range("H116").Value = 5993.6
$ws.Range("I116").Value = 5992.25
$ws.Range("K116").Value = 5992.25
$ws.Range("M116").Value = -2550.25

$ws.Range("H122").Value = 1290.5927
$ws.Range("I122").Value = 1313.88
$ws.Range("K122").Value = 3941.64
$ws.Range("M122").Value = -1491.64

$ws.Range("H129").Value = 596.8570999999999
$ws.Range("I129").Value = 596.8570999999999
$ws.Range("K129").Value = 1790.5713
$ws.Range("M129").Value = 3209.4287

$ws.Range("H137").Value = 2111.111
$ws.Range("J137").Value = 2166.3333
$ws.Range("L137").Value = 6498.999899999999
$ws.Range("N137").Value = -11598.9999

$ws.Range("H138").Value = 2394.4915
$ws.Range("J138").Value = 2839.9355
$ws.Range("L138").Value = 8519.806500000001
$ws.Range("N138").Value = -18799.8065

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("I2").Value = 25641634
$ws.Range("K2").Value = 25641634
$ws.Range("M2").Value = -25641521

$ws.Range("H74").Value = 522.4815
$ws.Range("I74").Value = 522.4815
$ws.Range("K74").Value = 522.4815
$ws.Range("M74").Value = 351.5185

$ws.Range("H77").Value = 522.4815
$ws.Range("I77").Value = 522.4815
$ws.Range("K77").Value = 2612.4075
$ws.Range("M77").Value = 1755.5925

$ws.Range("I116").Value = 25641634
$ws.Range("K116").Value = 25641634
$ws.Range("M116").Value = -25639340

$ws.Range("H132").Value = 35004
$ws.Range("I132").Value = 12
$ws.Range("J132").Value = 52500
$ws.Range("K132").Value = 36
$ws.Range("L132").Value = 157500
$ws.Range("M132").Value = 2494
$ws.Range("N132").Value = -162560

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("I3").Value = 25641634
$ws.Range("K3").Value = 25641634
$ws.Range("M3").Value = -25641520

$ws.Range("H99").Value = 1842.9412
$ws.Range("I99").Value = 1486.9231
$ws.Range("K99").Value = 1486.9231
$ws.Range("M99").Value = 11.07690000000002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 694.5
$ws.Range("I22").Value = 694.5
$ws.Range("K22").Value = 694.5
$ws.Range("M22").Value = -344.5

$ws.Range("H31").Value = 3264.9644
$ws.Range("I31").Value = 2216.2856
$ws.Range("K31").Value = 2216.2856
$ws.Range("M31").Value = -1921.2856

$ws.Range("H34").Value = 3264.9644
$ws.Range("I34").Value = 2216.2856
$ws.Range("K34").Value = 2216.2856
$ws.Range("M34").Value = -2014.2856

$ws.Range("H86").Value = 8498.5
$ws.Range("J86").Value = 9997
$ws.Range("L86").Value = 9997
$ws.Range("N86").Value = -12243

$ws.Range("H89").Value = 8498.5
$ws.Range("J89").Value = 9997
$ws.Range("L89").Value = 49985
$ws.Range("N89").Value = -61217

$ws.Range("H99").Value = 8688.324000000001
$ws.Range("I99").Value = 5033.1816
$ws.Range("K99").Value = 5033.1816
$ws.Range("M99").Value = -3535.1816

$ws.Range("H126").Value = 8688.324000000001
$ws.Range("I126").Value = 5033.1816
$ws.Range("K126").Value = 15099.5448
$ws.Range("M126").Value = -12629.5448

$ws.Range("H132").Value = 2383.5
$ws.Range("I132").Value = 1997.2727
$ws.Range("J132").Value = 3799.6667
$ws.Range("K132").Value = 5991.8181
$ws.Range("L132").Value = 11399.0001
$ws.Range("M132").Value = -3461.8181
$ws.Range("N132").Value = -16459.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 67137.47
$ws.Range("I2").Value = 125077.875
$ws.Range("J2").Value = 919.8570999999999
$ws.Range("K2").Value = 750467.25
$ws.Range("L2").Value = 5519.142599999999
$ws.Range("M2").Value = -750354.25
$ws.Range("N2").Value = -5745.142599999999

$ws.Range("H5").Value = 1313.3334
$ws.Range("I5").Value = 1470
$ws.Range("J5").Value = 1000
$ws.Range("K5").Value = 4410
$ws.Range("L5").Value = 3000
$ws.Range("M5").Value = -4298
$ws.Range("N5").Value = -3224

$ws.Range("H56").Value = 12316.9
$ws.Range("I56").Value = 12316.9
$ws.Range("K56").Value = 12316.9
$ws.Range("M56").Value = -11786.9

$ws.Range("H93").Value = 4000
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 4000
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 12000
$ws.Range("M93").ClearContents()
$ws.Range("N93").Value = -15744

$ws.Range("H111").Value = 3864
$ws.Range("J111").Value = 3864
$ws.Range("L111").Value = 11592
$ws.Range("N111").Value = -17726

$ws.Range("H134").Value = 2030
$ws.Range("I134").Value = 2030
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 6090
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -1020
$ws.Range("N134").ClearContents()

$ws.Range("H135").Value = 1313.3334
$ws.Range("I135").Value = 1470
$ws.Range("J135").Value = 1000
$ws.Range("K135").Value = 13230
$ws.Range("L135").Value = 9000
$ws.Range("M135").Value = -10695
$ws.Range("N135").Value = -14070

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4800
$ws.Range("I70").Value = 4800
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 4800
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -4530
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 4800
$ws.Range("I73").Value = 4800
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 4800
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -3864
$ws.Range("N73").ClearContents()

$ws.Range("H132").Value = 1937.375
$ws.Range("I132").Value = 1937.375
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5812.125
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -3282.125
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 999999
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 999999
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 999999
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -1000375

$ws.Range("H55").Value = 533.1818
$ws.Range("I55").Value = 283.5
$ws.Range("K55").Value = 283.5
$ws.Range("M55").Value = -110.5

$ws.Range("H132").Value = 3333.3333
$ws.Range("I132").Value = 1000
$ws.Range("J132").Value = 4500
$ws.Range("K132").Value = 3000
$ws.Range("L132").Value = 13500
$ws.Range("M132").Value = -470
$ws.Range("N132").Value = -18560

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2539.111
$ws.Range("I132").Value = 1880.3334
$ws.Range("K132").Value = 5641.0002
$ws.Range("M132").Value = -3111.0002

$ws.Range("H136").Value = 1804.3077
$ws.Range("I136").Value = 1145.6
$ws.Range("K136").Value = 3436.8
$ws.Range("M136").Value = -886.7999999999997
